# Clean up stray leading/trailing whitespace in the "SSPLact" / "SSPPreg"
# age-group labels (rows 45-46, column A of Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A45").Value = "SSPLact"
$ws.Range("A46").Value = "SSPPreg"

# Leave the selection on the last touched cell, matching the saved file.
$ws.Range("A46").Select()
